# Restructure sheet "Q145" to add "All"/"Men"/"Women" breakdowns, per commit:
# "Added support to differentiation by gender"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the existing ("All") table down by one row and label it ---
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "All"

# --- 2. "Men" table starting at row 18 ---
$ws.Range("A18").Value = "Men"

$menHeader = New-Object 'object[,]' 1,3
$menHeader[0,0] = 'response'
$menHeader[0,1] = 'Uganda (% of respondents)'
$menHeader[0,2] = 'Kenya (% of respondents)'
$ws.Range("B19:D19").Value = $menHeader

$menData = New-Object 'object[,]' 14,4
$menData[0,0] = 1
$menData[0,1] = 'Being a member of a local group (women''s groups etc)'
$menData[0,2] = 0.14070351758794
$menData[0,3] = 0.00641025641025641
$menData[1,0] = 2
$menData[1,1] = 'Centres of Excellence'
$menData[1,2] = 0.035175879396984903
$menData[1,3] = 0.019230769230769201
$menData[2,0] = 3
$menData[2,1] = 'Extension workers from your cooperative'
$menData[2,2] = 0.48241206030150802
$menData[2,3] = 0.262820512820513
$menData[3,0] = 4
$menData[3,1] = 'Guidance from farmer promoters'
$menData[3,2] = 0.221105527638191
$menData[3,3] = 0.32692307692307698
$menData[4,0] = 5
$menData[4,1] = 'Internet'
$menData[4,2] = 0.0050251256281407001
$menData[4,3] = 0.025641025641025599
$menData[5,0] = 6
$menData[5,1] = 'Mobile phones'
$menData[5,2] = 0.0552763819095477
$menData[5,3] = 0.083333333333333301
$menData[6,0] = 7
$menData[6,1] = 'My own experimentation'
$menData[6,2] = 0.42211055276381898
$menData[6,3] = 0.30128205128205099
$menData[7,0] = 8
$menData[7,1] = 'Newspapers'
$menData[7,2] = 0.050251256281407003
$menData[7,3] = 0.108974358974359
$menData[8,0] = 9
$menData[8,1] = 'Radio stations'
$menData[8,2] = 0.34673366834170899
$menData[8,3] = 0.56410256410256399
$menData[9,0] = 10
$menData[9,1] = 'Sharing with other farmers/friends/family'
$menData[9,2] = 0.42211055276381898
$menData[9,3] = 0.35897435897435898
$menData[10,0] = 11
$menData[10,1] = 'Television'
$menData[10,2] = 0.045226130653266298
$menData[10,3] = 0.121794871794872
$menData[11,0] = 12
$menData[11,1] = 'Training materials / leaflets'
$menData[11,2] = 0.035175879396984903
$menData[11,3] = 0.089743589743589702
$menData[12,0] = 13
$menData[12,1] = 'Training sessions/workshops (including FFSs)'
$menData[12,2] = 0.32663316582914598
$menData[12,3] = 0.269230769230769
$menData[13,0] = 14
$menData[13,1] = 'WeFarm'
$menData[13,2] = 0.0100502512562814
$menData[13,3] = 0.0448717948717949
$ws.Range("A20:D33").Value = $menData

# --- 3. "Women" table starting at row 35 ---
$ws.Range("A35").Value = "Women"

$womenHeader = New-Object 'object[,]' 1,3
$womenHeader[0,0] = 'response'
$womenHeader[0,1] = 'Uganda (% of respondents)'
$womenHeader[0,2] = 'Kenya (% of respondents)'
$ws.Range("B36:D36").Value = $womenHeader

$womenData = New-Object 'object[,]' 14,4
$womenData[0,0] = 1
$womenData[0,1] = 'Being a member of a local group (women''s groups etc)'
$womenData[0,2] = 0.160377358490566
$womenData[0,3] = 0.0204081632653061
$womenData[1,0] = 2
$womenData[1,1] = 'Centres of Excellence'
$womenData[1,2] = 0.066037735849056603
$womenData[1,3] = 0.0408163265306122
$womenData[2,0] = 3
$womenData[2,1] = 'Extension workers from your cooperative'
$womenData[2,2] = 0.52830188679245305
$womenData[2,3] = 0.122448979591837
$womenData[3,0] = 4
$womenData[3,1] = 'Guidance from farmer promoters'
$womenData[3,2] = 0.34905660377358499
$womenData[3,3] = 0.38775510204081598
$womenData[4,0] = 5
$womenData[4,1] = 'Internet'
$womenData[4,2] = 0
$womenData[4,3] = 0
$womenData[5,0] = 6
$womenData[5,1] = 'Mobile phones'
$womenData[5,2] = 0.0188679245283019
$womenData[5,3] = 0.122448979591837
$womenData[6,0] = 7
$womenData[6,1] = 'My own experimentation'
$womenData[6,2] = 0.56603773584905703
$womenData[6,3] = 0.34693877551020402
$womenData[7,0] = 8
$womenData[7,1] = 'Newspapers'
$womenData[7,2] = 0.0188679245283019
$womenData[7,3] = 0.081632653061224497
$womenData[8,0] = 9
$womenData[8,1] = 'Radio stations'
$womenData[8,2] = 0.41509433962264197
$womenData[8,3] = 0.59183673469387799
$womenData[9,0] = 10
$womenData[9,1] = 'Sharing with other farmers/friends/family'
$womenData[9,2] = 0.57547169811320797
$womenData[9,3] = 0.42857142857142899
$womenData[10,0] = 11
$womenData[10,1] = 'Television'
$womenData[10,2] = 0.0094339622641509396
$womenData[10,3] = 0.14285714285714299
$womenData[11,0] = 12
$womenData[11,1] = 'Training materials / leaflets'
$womenData[11,2] = 0.056603773584905703
$womenData[11,3] = 0.0204081632653061
$womenData[12,0] = 13
$womenData[12,1] = 'Training sessions/workshops (including FFSs)'
$womenData[12,2] = 0.320754716981132
$womenData[12,3] = 0.16326530612244899
$womenData[13,0] = 14
$womenData[13,1] = 'WeFarm'
$womenData[13,2] = 0
$womenData[13,3] = 0.0408163265306122
$ws.Range("A37:D50").Value = $womenData


# --- 4. Re-point the chart's series at the shifted "All" range ---
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection()
$s1 = $series.Item(1)
$s2 = $series.Item(2)
$s1.Formula = "=SERIES('Q145'!`$C`$2,'Q145'!`$B`$3:`$B`$16,'Q145'!`$C`$3:`$C`$16,1)"
$s2.Formula = "=SERIES('Q145'!`$D`$2,'Q145'!`$B`$3:`$B`$16,'Q145'!`$D`$3:`$D`$16,2)"

# --- 5. Move/resize the chart to its new anchor ---
$chartObj.Left = 327.1875
$chartObj.Top = 22
$chartObj.Width = 508.9375
$chartObj.Height = 363.35291338582675

# --- 6. Restore the reported selection ---
$ws.Range("E10").Select()
